$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StringLocalizations_BasicText")
[void]$ws.Activate()

# New "tips" rows appended after the existing data (rows 106-110), matching
# the key/value layout used throughout the sheet (A=key, B=value,
# C/D/E=untranslated placeholder "XXXX"). The string cells are written in
# the same first-seen order as the authored edit so the shared-string table
# comes out in the same sequence.
$ws.Range("A106").Value = "TIPS_SATISFACTION"
$ws.Range("B106").Value = "{0}% Citizen Satisfaction*n*Keep the badge out of the red area by making good choices"
$ws.Range("B108").Value = "{0} Active Cases*n*Close cases by making good choices to make a big difference to citizen satisfaction"
$ws.Range("A107").Value = "TIPS_INCIDENT"
$ws.Range("A109").Value = "TIPS_OFFICER"
$ws.Range("A110").Value = "TIPS_OFFICERS"
$ws.Range("A108").Value = "TIPS_INCIDENTS"
$ws.Range("B107").Value = "{0} Active Case*n*Close cases by making good choices to make a big difference to citizen satisfaction"
$ws.Range("B109").Value = "{0} Available Officer*n*A limited number available so be careful, you never know when a severe incident may occur"
$ws.Range("B110").Value = "{0} Available Officers*n*A limited number available so be careful, you never know when a severe incident may occur"

for ($rowNum = 106; $rowNum -le 110; $rowNum++) {
    $ws.Range("C$rowNum").Value = "XXXX"
    $ws.Range("D$rowNum").Value = "XXXX"
    $ws.Range("E$rowNum").Value = "XXXX"

    $rowRange = $ws.Range("A${rowNum}:E${rowNum}")
    $rowRange.VerticalAlignment = -4160   # xlTop
    $rowRange.WrapText = $false

    $ws.Range("B$rowNum").WrapText = $true

    $ws.Rows.Item($rowNum).RowHeight = 28.8
}

# Reposition the view / selection like the authored edit did.
$excel.ActiveWindow.ScrollRow = 94
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B108").Select()

Write-Output "Added tips rows 106-110"
